$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates) in the new rows should pick up the same date-formatted
# style already used for the rest of the table (s=13), so copy the format
# from an existing dated cell before filling in the values.
$ws.Range("A3").Copy()
$ws.Range("A17:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 17 - 3/18/2017, 4 hours
$ws.Range("A17").Value = "3/18/2017"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = "Learning JSP basics"

# Row 18 - 3/31/2017, 2 hours
$ws.Range("A18").Value = "3/31/2017"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Setting up development environment and server"

# Row 19 - 4/23/2017, 6 hours
$ws.Range("A19").Value = "4/23/2017"
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = "Updated user login and create user page"

# Row 22 - 4/27/2017, 6 hours (filled in out of order, ahead of rows 20-21)
$ws.Range("A22").Value = "4/27/2017"
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = "Added CSS and changed some html design"

# Row 20 - 4/25/2017, 6 hours
$ws.Range("A20").Value = "4/25/2017"
$ws.Range("B20").Value = 6
$ws.Range("C20").Value = "Added servlets and updated xml"

# Row 21 - 4/26/2017, 2 hours
$ws.Range("A21").Value = "4/26/2017"
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "Prepared for presentation"

# B19 picked up a slightly different border treatment (no top/bottom rule)
# versus the rest of the Hours column - drop the horizontal edges so it
# matches.
$rng = $ws.Range("B19")
$rng.Borders.Item(9).LineStyle = -4142
$rng.Borders.Item(8).LineStyle = -4142

# Leave the selection where the user finished editing.
$null = $ws.Range("C22").Select()

Write-Output "done"
